$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Virtual Machines": the Private DNS hostnames for the heider-temp-vm
# and custon-dns-temp rows are no longer valid, so they become "N/A".
# ---------------------------------------------------------------------------
$wsVm = $wb.Worksheets.Item("Virtual Machines")
$wsVm.Cells.Item(4, 5).Value = "N/A"
$wsVm.Cells.Item(5, 5).Value = "N/A"

# ---------------------------------------------------------------------------
# Sheet "AKS": replace the single-cluster service listing with the gathered
# IPs for the two real AKS clusters (aks-heider-temp / aks-test2-heider).
# ---------------------------------------------------------------------------
$wsAks = $wb.Worksheets.Item("AKS")

$aksRows = @(
    @("Resource Group", "AKS Server", "Namespace", "Service", "Service IP"),
    @("MC_aks-heider_aks-heider-temp_eastus", "aks-heider-temp", "calico-system", "calico-kube-controllers-metrics", "None"),
    @("MC_aks-heider_aks-heider-temp_eastus", "aks-heider-temp", "calico-system", "calico-typha", "10.0.226.138"),
    @("MC_aks-heider_aks-heider-temp_eastus", "aks-heider-temp", "default", "kubernetes", "10.0.0.1"),
    @("MC_aks-heider_aks-heider-temp_eastus", "aks-heider-temp", "kube-system", "kube-dns", "10.0.0.10"),
    @("MC_aks-heider_aks-heider-temp_eastus", "aks-heider-temp", "kube-system", "metrics-server", "10.0.2.10"),
    @("MC_aks-heider_aks-test2-heider_westus", "aks-test2-heider", "calico-system", "calico-kube-controllers-metrics", "None"),
    @("MC_aks-heider_aks-test2-heider_westus", "aks-test2-heider", "calico-system", "calico-typha", "10.0.226.138"),
    @("MC_aks-heider_aks-test2-heider_westus", "aks-test2-heider", "default", "kubernetes", "10.0.0.1"),
    @("MC_aks-heider_aks-test2-heider_westus", "aks-test2-heider", "kube-system", "kube-dns", "10.0.0.10"),
    @("MC_aks-heider_aks-test2-heider_westus", "aks-test2-heider", "kube-system", "metrics-server", "10.0.2.10")
)

for ($r = 1; $r -le $aksRows.Length; $r++) {
    $rowVals = $aksRows[$r - 1]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $wsAks.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Leave the selection on A1 (old sheet had a full-column selection lingering
# from before; put the cursor back to the top-left like a freshly edited
# sheet would have).
$wsAks.Range("A1").Select()

# Resize the helper columns so the new, longer resource-group / server names
# are not truncated (mirrors the width bump seen after AutoFit in Excel).
$wsAks.Columns.Item(1).ColumnWidth = 35.1796875
$wsAks.Columns.Item(2).ColumnWidth = 14.6328125
$wsAks.Columns.Item(3).ColumnWidth = 11.90625
